$wb = $excel.ActiveWorkbook

# Resize/reposition the workbook window (mirrors the saved bookViews/workbookView).
$win = $excel.ActiveWindow
$win.Left = 4080
$win.Top = 0
$win.Width = 24720
$win.Height = 16740

# "openbis-metadata" sheet: fill in the Strain and Experiment values.
$ws = $wb.Worksheets.Item("openbis-metadata")
$ws.Activate()

# Set B3 (Strain) first, then B2 (Experiment) so new shared strings are
# interned in the same order as the target workbook (MGP9 before
# /TEST/TEST/TEST).
$ws.Range("B3").Value = "MGP9"
$ws.Range("B2").Value = "/TEST/TEST/TEST"

$ws.Range("B3").Select()
